$d = $word.ActiveDocument

# Remove the first paragraph (date + name/title line) entirely, including
# its paragraph mark, so the following "Data/AI Engineer..." paragraph
# moves up and becomes the first paragraph.
$p1 = $d.Paragraphs(1)
$p1.Range.Delete()

# The paragraph that moved up keeps its own (BodyText) style after the
# merge, so restore the FirstParagraph style that the original opening
# paragraph used.
$d.Paragraphs(1).Style = "First Paragraph"
